# Generate Report for handoff
# Replaces the stale "cf29bef5-..." handoff file references with the new
# "6d7d606d-..." ones, updates the handoff/target timestamps, and removes
# the row for the file that previously failed transform (which no longer
# appears in the regenerated report), so each sheet shrinks from 4 data
# rows to 3.

$wb = $excel.ActiveWorkbook

$newMd  = "6d7d606d-b4da-46d5-aa6f-091d4023c639.md"

$newMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/16ee7691d739ee6da4410eca9f4b6ae20088658b/e2e/6d7d606d-b4da-46d5-aa6f-091d4023c639.md"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/16ee7691d739ee6da4410eca9f4b6ae20088658b/.localization-config"

$newZhXlf = "6d7d606d-b4da-46d5-aa6f-091d4023c639.a3a532de84f86b30cbad759a6a4157e18d00e82c.zh-cn.xlf"
$newDeXlf = "6d7d606d-b4da-46d5-aa6f-091d4023c639.a3a532de84f86b30cbad759a6a4157e18d00e82c.de-de.xlf"

$newZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c5f29c90693315bdf3f15ada963ebcd4a32263a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/6d7d606d-b4da-46d5-aa6f-091d4023c639.a3a532de84f86b30cbad759a6a4157e18d00e82c.zh-cn.xlf"
$newDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4f5bf454ef659b777d0b7a71fe5e0c7905f23e1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/6d7d606d-b4da-46d5-aa6f-091d4023c639.a3a532de84f86b30cbad759a6a4157e18d00e82c.de-de.xlf"

$zhDatetime = "2016-01-13 11:33:46"
$deDatetime = "2016-01-13 11:33:59"

$cfg    = ".localization-config"
$notLoc = "Not to be localized"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() removes every hyperlink on the sheet, so clear them
# all up-front and rebuild the ones that remain afterwards.
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = $newMd
$ws1.Range("A3").Value = $cfg
$ws1.Range("B3").Value = $notLoc
$ws1.Range("C3").Value = $notLoc

# Drop the old row 3 (the "Handoff transform failed" entry) - row 4
# (.localization-config) shifts up into its place, so after updating row 3
# above to hold the .localization-config data we delete what is now the
# spare trailing row.
$ws1.Rows.Item(4).Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfg) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = $newMd
$ws2.Range("C2").Value = $newZhXlf
$ws2.Range("D2").Value = $zhDatetime

$ws2.Range("A3").Value = $cfg
$ws2.Range("B3").Value = $notLoc

$ws2.Rows.Item(4).Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $newZhXlfUrl, [Type]::Missing, [Type]::Missing, $newZhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfg) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = $newMd
$ws3.Range("C2").Value = $newDeXlf
$ws3.Range("D2").Value = $deDatetime

$ws3.Range("A3").Value = $cfg
$ws3.Range("B3").Value = $notLoc

$ws3.Rows.Item(4).Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $newDeXlfUrl, [Type]::Missing, [Type]::Missing, $newDeXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfg) | Out-Null
